$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "AFMH Farm Coomboona",
    "Barwon Heads Hotel Barwon Heads",
    "Confirmed Omicron Sircuit Bar Fitzroy",
    "Confirmed Omicron Variant The Peel Hotel Collingwood",
    "Costa Mushroom Farm Mernda",
    "Flanagans Border Inn Bacchus Marsh",
    "Hurley's Bar and Bistro Hotel Kyabram",
    "Melbourne Cricket Ground (MCG)",
    "Monash Health Dandenong Hospital Dandenong",
    "Ms Collins Melbourne",
    "Northern Health Northern Hospital Epping",
    "Rupert On Rupert Collingwood",
    "St Vincents Hospital Melbourne Emergency Department Fitzroy",
    "St. Vincent's Hospital Melbourne Fitzroy",
    "The Deck Shepparton",
    "Werribee Mercy Hospital Emergency Department",
    "Western Health Sunshine Hospital Emergency Department St Albans",
    "Wilson's Nightclub Horsham"
)

$values = @(18, 29, 17, 10, 10, 16, 10, 65, 16, 24, 11, 11, 23, 12, 12, 10, 10, 13)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
